$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (shift left by one quarter, add new quarter in column M) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "3 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "6 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "9 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "12 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "3 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "6 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "9 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "12 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "3 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1400-12-02 (9)"
$ws.Range("E9").Value = "1401-01-30 (2)"
$ws.Range("F9").Value = "1401-05-31 (6)"
$ws.Range("G9").Value = "1401-07-30 (2)"
$ws.Range("H9").Value = "1401-12-12 (8)"
$ws.Range("I9").Value = "1402-01-30 (2)"
$ws.Range("J9").Value = "1401-05-31 (2)"
$ws.Range("K9").Value = "1401-07-30"
$ws.Range("L9").Value = "1402-01-30 (3)"
$ws.Range("M9").Value = "1402-01-30"

# --- Numeric data rows ---
# Row 11
$ws.Range("D11").Value = 53901
$ws.Range("E11").Value = 16941
$ws.Range("F11").Value = 41285
$ws.Range("G11").Value = 58548
$ws.Range("H11").Value = 77688
$ws.Range("I11").Value = 18469
$ws.Range("J11").Value = 42378
$ws.Range("K11").Value = 63613
$ws.Range("L11").Value = 86789
$ws.Range("M11").Value = 17169

# Row 12
$ws.Range("D12").Value = -40789
$ws.Range("E12").Value = -11101
$ws.Range("F12").Value = -27928
$ws.Range("G12").Value = -41417
$ws.Range("H12").Value = -57098
$ws.Range("I12").Value = -14485
$ws.Range("J12").Value = -36069
$ws.Range("K12").Value = -53792
$ws.Range("L12").Value = -71329
$ws.Range("M12").Value = -12537

# Row 13
$ws.Range("D13").Value = 13112
$ws.Range("E13").Value = 5840
$ws.Range("F13").Value = 13356
$ws.Range("G13").Value = 17132
$ws.Range("H13").Value = 20590
$ws.Range("I13").Value = 3984
$ws.Range("J13").Value = 6308
$ws.Range("K13").Value = 9821
$ws.Range("L13").Value = 15460
$ws.Range("M13").Value = 4632

# Row 14
$ws.Range("D14").Value = -3110
$ws.Range("E14").Value = -548
$ws.Range("F14").Value = -2030
$ws.Range("G14").Value = -2777
$ws.Range("H14").Value = -4127
$ws.Range("I14").Value = -927
$ws.Range("J14").Value = -2642
$ws.Range("K14").Value = -3918
$ws.Range("L14").Value = -5525
$ws.Range("M14").Value = -966

# Row 16
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 259
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 541
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = 291
$ws.Range("K16").Value = 346
$ws.Range("L16").Value = 1116
$ws.Range("M16").Value = 152

# Row 17
$ws.Range("D17").Value = 10002
$ws.Range("E17").Value = 5293
$ws.Range("F17").Value = 11585
$ws.Range("G17").Value = 14356
$ws.Range("H17").Value = 17004
$ws.Range("I17").Value = 3175
$ws.Range("J17").Value = 3957
$ws.Range("K17").Value = 6250
$ws.Range("L17").Value = 11051
$ws.Range("M17").Value = 3818

# Row 18
$ws.Range("D18").Value = -1935
$ws.Range("E18").Value = -894
$ws.Range("F18").Value = -1737
$ws.Range("G18").Value = -2445
$ws.Range("H18").Value = -3759
$ws.Range("I18").Value = -1255
$ws.Range("J18").Value = -2596
$ws.Range("K18").Value = -3439
$ws.Range("L18").Value = -4769
$ws.Range("M18").Value = -1056

# Row 19
$ws.Range("D19").Value = 1851
$ws.Range("E19").Value = 204
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 464
$ws.Range("H19").Value = 66
$ws.Range("I19").Value = 129
$ws.Range("J19").Value = 272
$ws.Range("K19").Value = 507
$ws.Range("L19").Value = 178
$ws.Range("M19").Value = 6

# Row 20
$ws.Range("D20").Value = 9918
$ws.Range("E20").Value = 4603
$ws.Range("F20").Value = 9905
$ws.Range("G20").Value = 12375
$ws.Range("H20").Value = 13310
$ws.Range("I20").Value = 2050
$ws.Range("J20").Value = 1633
$ws.Range("K20").Value = 3317
$ws.Range("L20").Value = 6460
$ws.Range("M20").Value = 2768

# Row 21
$ws.Range("D21").Value = -3785
$ws.Range("E21").Value = -1017
$ws.Range("F21").Value = -1546
$ws.Range("G21").Value = -2039
$ws.Range("H21").Value = -1952
$ws.Range("I21").Value = -461
$ws.Range("J21").Value = -309
$ws.Range("K21").Value = -601
$ws.Range("L21").Value = -1205
$ws.Range("M21").Value = -550

# Row 22
$ws.Range("D22").Value = 6133
$ws.Range("E22").Value = 3587
$ws.Range("F22").Value = 8359
$ws.Range("G22").Value = 10336
$ws.Range("H22").Value = 11359
$ws.Range("I22").Value = 1589
$ws.Range("J22").Value = 1324
$ws.Range("K22").Value = 2715
$ws.Range("L22").Value = 5254
$ws.Range("M22").Value = 2218

# Row 24
$ws.Range("D24").Value = 6133
$ws.Range("E24").Value = 3587
$ws.Range("F24").Value = 8359
$ws.Range("G24").Value = 10336
$ws.Range("H24").Value = 11359
$ws.Range("I24").Value = 1589
$ws.Range("J24").Value = 1324
$ws.Range("K24").Value = 2715
$ws.Range("L24").Value = 5254
$ws.Range("M24").Value = 2218

# Row 26
$ws.Range("D26").Value = 10943
$ws.Range("E26").Value = 8926
$ws.Range("F26").Value = 9128
$ws.Range("G26").Value = 8842
$ws.Range("H26").Value = 8511
$ws.Range("I26").Value = 7915
$ws.Range("J26").Value = 7670
$ws.Range("K26").Value = 14901
$ws.Range("L26").Value = 14204
$ws.Range("M26").Value = 9592

